$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated team-specific time-transition matrix values (Miami (FL)_A)
$ws.Range("B2").Value = 0.1476793248945148
$ws.Range("C2").Value = 0.6540084388185654
$ws.Range("J2").Value = 0.01265822784810127
$ws.Range("P2").Value = 0.09282700421940929
$ws.Range("S2").Value = 0.09282700421940929
$ws.Range("C3").Value = 0.01875
$ws.Range("J3").Value = 0.0375
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.19375
$ws.Range("P4").Value = 0.7027027027027027
$ws.Range("S4").Value = 0.2972972972972973
$ws.Range("B6").Value = 0.04716981132075472
$ws.Range("D6").Value = 0.009433962264150943
$ws.Range("F6").Value = 0.03773584905660377
$ws.Range("J6").Value = 0.2311320754716981
$ws.Range("O6").Value = 0.009433962264150943
$ws.Range("Q6").Value = 0.1981132075471698
$ws.Range("R6").Value = 0.1037735849056604
$ws.Range("S6").Value = 0.3632075471698113
$ws.Range("B7").Value = 0.07194244604316546
$ws.Range("D7").Value = 0.007194244604316547
$ws.Range("F7").Value = 0.07194244604316546
$ws.Range("J7").Value = 0.1151079136690648
$ws.Range("O7").Value = 0.02877697841726619
$ws.Range("Q7").Value = 0.1798561151079137
$ws.Range("R7").Value = 0.1223021582733813
$ws.Range("S7").Value = 0.4028776978417266
$ws.Range("B8").Value = 0.07838479809976247
$ws.Range("D8").Value = 0.01425178147268409
$ws.Range("F8").Value = 0.05225653206650831
$ws.Range("J8").Value = 0.1377672209026128
$ws.Range("O8").Value = 0.01187648456057007
$ws.Range("Q8").Value = 0.1900237529691211
$ws.Range("R8").Value = 0.1140142517814727
$ws.Range("S8").Value = 0.4014251781472684
$ws.Range("B9").Value = 0.1056910569105691
$ws.Range("D9").Value = 0.008130081300813009
$ws.Range("F9").Value = 0.05691056910569105
$ws.Range("J9").Value = 0.1382113821138211
$ws.Range("O9").Value = 0.01219512195121951
$ws.Range("Q9").Value = 0.1463414634146341
$ws.Range("R9").Value = 0.1219512195121951
$ws.Range("S9").Value = 0.4105691056910569
$ws.Range("B10").Value = 0.08664259927797834
$ws.Range("D10").Value = 0.01949458483754513
$ws.Range("F10").Value = 0.07075812274368232
$ws.Range("J10").Value = 0.1335740072202166
$ws.Range("O10").Value = 0.01444043321299639
$ws.Range("Q10").Value = 0.2194945848375451
$ws.Range("R10").Value = 0.1090252707581227
$ws.Range("S10").Value = 0.3465703971119133
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.06190476190476191
$ws.Range("K11").Value = 0.2095238095238095
$ws.Range("L11").Value = 0.5523809523809524
$ws.Range("S11").Value = 0.009523809523809525
$ws.Range("G12").Value = 0.6557377049180327
$ws.Range("J12").Value = 0.2377049180327869
$ws.Range("K12").Value = 0.01639344262295082
$ws.Range("L12").Value = 0.05737704918032787
$ws.Range("S12").Value = 0.03278688524590164
$ws.Range("G13").Value = 0.7297297297297297
$ws.Range("J13").Value = 0.2162162162162162
$ws.Range("S13").Value = 0.05405405405405406
$ws.Range("F15").Value = 0.0198019801980198
$ws.Range("H15").Value = 0.1732673267326733
$ws.Range("I15").Value = 0.07425742574257425
$ws.Range("J15").Value = 0.3316831683168317
$ws.Range("K15").Value = 0.06435643564356436
$ws.Range("M15").Value = 0.0198019801980198
$ws.Range("N15").Value = 0.004950495049504951
$ws.Range("O15").Value = 0.0891089108910891
$ws.Range("S15").Value = 0.2227722772277228
$ws.Range("F16").Value = 0.01257861635220126
$ws.Range("H16").Value = 0.1509433962264151
$ws.Range("I16").Value = 0.1069182389937107
$ws.Range("J16").Value = 0.4905660377358491
$ws.Range("K16").Value = 0.08176100628930817
$ws.Range("M16").Value = 0.01886792452830189
$ws.Range("N16").Value = 0.006289308176100629
$ws.Range("O16").Value = 0.06289308176100629
$ws.Range("S16").Value = 0.06918238993710692
$ws.Range("F17").Value = 0.01652892561983471
$ws.Range("H17").Value = 0.1983471074380165
$ws.Range("I17").Value = 0.08677685950413223
$ws.Range("J17").Value = 0.4690082644628099
$ws.Range("K17").Value = 0.04958677685950413
$ws.Range("M17").Value = 0.01239669421487603
$ws.Range("O17").Value = 0.05785123966942149
$ws.Range("S17").Value = 0.109504132231405
$ws.Range("F18").Value = 0.01503759398496241
$ws.Range("H18").Value = 0.1390977443609022
$ws.Range("I18").Value = 0.09398496240601503
$ws.Range("J18").Value = 0.5150375939849624
$ws.Range("K18").Value = 0.04511278195488722
$ws.Range("M18").Value = 0.02255639097744361
$ws.Range("O18").Value = 0.05263157894736842
$ws.Range("S18").Value = 0.1165413533834586
$ws.Range("F19").Value = 0.01221498371335505
$ws.Range("H19").Value = 0.1848534201954397
$ws.Range("I19").Value = 0.1205211726384365
$ws.Range("J19").Value = 0.4079804560260586
$ws.Range("K19").Value = 0.07736156351791532
$ws.Range("M19").Value = 0.01465798045602606
$ws.Range("O19").Value = 0.06107491856677524
$ws.Range("S19").Value = 0.1213355048859935
